# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and swap the Hedera/WhiteBITCoin rows (46 <-> 47) to match the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.497.05'
$ws.Range("E2").Value = '  -3.53%  '

$ws.Range("D3").Value = '2.648.98'
$ws.Range("E3").Value = '  -1.97%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = "'522.01"
$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").Value = "'144.63"
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = '  -1.70%  '

$ws.Range("D9").Value = "'6.71"
$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("E10").Value = '  -2.91%  '

$ws.Range("E11").Value = '  -0.77%  '

$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("D13").Value = '3.114.42'
$ws.Range("E13").Value = '  -2.11%  '

$ws.Range("D14").Value = '58.500.54'
$ws.Range("E14").Value = '  -3.84%  '

$ws.Range("D15").Value = "'20.89"
$ws.Range("E15").Value = '  -2.03%  '

$ws.Range("D16").Value = "'0.0000137"
$ws.Range("E16").Value = '  -1.19%  '

$ws.Range("D17").Value = '2.658.07'
$ws.Range("E17").Value = '  -8.10%  '

$ws.Range("D18").Value = "'338.57"
$ws.Range("E18").Value = '  -2.73%  '

$ws.Range("D19").Value = "'4.42"
$ws.Range("E19").Value = '  -2.15%  '

$ws.Range("D20").Value = "'10.49"
$ws.Range("E20").Value = '  -1.16%  '

$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = '  -1.58%  '

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").Value = "'64.52"
$ws.Range("E23").Value = '  +1.11%  '

$ws.Range("D24").Value = "'0.425"
$ws.Range("E24").Value = '  +1.06%  '

$ws.Range("E25").Value = '  -1.94%  '

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = '  +0.59%  '

$ws.Range("D27").Value = '0.0₃0799'
$ws.Range("E27").Value = '  -2.41%  '

$ws.Range("D28").Value = "'7.13"
$ws.Range("E28").Value = '  -2.48%  '

$ws.Range("D29").Value = "'6.64"
$ws.Range("E29").Value = '  -1.91%  '

$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("E31").Value = '  -1.03%  '

$ws.Range("D32").Value = "'152.68"
$ws.Range("E32").Value = '  +1.57%  '

$ws.Range("D33").Value = "'18.87"
$ws.Range("E33").Value = '  -1.43%  '

$ws.Range("D34").Value = "'4.15"
$ws.Range("E34").Value = '  -2.16%  '

$ws.Range("E35").Value = '  -3.73%  '

$ws.Range("D36").Value = "'0.910"
$ws.Range("E36").Value = '  -3.72%  '

$ws.Range("D37").Value = "'0.859"
$ws.Range("E37").Value = '  -2.01%  '

$ws.Range("D38").Value = "'36.68"
$ws.Range("E38").Value = '  -0.98%  '

$ws.Range("E39").Value = '  -4.65%  '

$ws.Range("E40").Value = '  -0.72%  '

$ws.Range("E41").Value = '  +0.40%  '

$ws.Range("D42").Value = "'0.608"
$ws.Range("E42").Value = '  -0.78%  '

$ws.Range("D43").Value = "'0.0970"
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").Value = "'270.46"
$ws.Range("E44").Value = '  -5.00%  '

$ws.Range("D45").Value = "'19.43"
$ws.Range("E45").Value = '  -3.65%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = "'0.0537"
$ws.Range("E46").Value = '  -0.43%  '

$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = "'10.64"
$ws.Range("E47").Value = '  +1.58%  '

$ws.Range("D48").Value = '2.047.93'
$ws.Range("E48").Value = '  -4.28%  '

$ws.Range("D49").Value = "'4.71"
$ws.Range("E49").Value = '  -1.88%  '

$ws.Range("E50").Value = '  -2.86%  '

$ws.Range("D51").Value = "'18.40"
$ws.Range("E51").Value = '  -5.00%  '

Write-Host "Updated 91 cells across 50 rows"
